$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41
$ws.Cells.Item(41, 2).Value = 6011491
$row41 = New-Object 'object[,]' 1,26
$row41[0,0] = "Sarpsborg"
$row41[0,1] = "HamKam"
$row41[0,2] = 2
$row41[0,3] = 3
$row41[0,4] = 2
$row41[0,5] = 3
$row41[0,6] = "A"
$row41[0,7] = 1.444
$row41[0,8] = 4.75
$row41[0,9] = 6.5
$row41[0,10] = 1.363
$row41[0,11] = 5.25
$row41[0,12] = 8
$row41[0,13] = -1.5
$row41[0,14] = 1.925
$row41[0,15] = 1.925
$row41[0,16] = 3.5
$row41[0,17] = 1.95
$row41[0,18] = 1.9
$row41[0,19] = -1
$row41[0,20] = -1
$row41[0,21] = 7
$row41[0,22] = -1
$row41[0,23] = 0.925
$row41[0,24] = 0.95
$row41[0,25] = -1
$ws.Range("E41:AD41").Value = $row41

# Row 44
$ws.Cells.Item(44, 2).Value = 6011402
$row44 = New-Object 'object[,]' 1,26
$row44[0,0] = "Aalesund"
$row44[0,1] = "Odd BK"
$row44[0,2] = 0
$row44[0,3] = 3
$row44[0,4] = 0
$row44[0,5] = 0
$row44[0,6] = "A"
$row44[0,7] = 2.3
$row44[0,8] = 3.3
$row44[0,9] = 3.1
$row44[0,10] = 2.3
$row44[0,11] = 3.3
$row44[0,12] = 3.3
$row44[0,13] = -0.25
$row44[0,14] = 1.975
$row44[0,15] = 1.875
$row44[0,16] = 2.5
$row44[0,17] = 1.975
$row44[0,18] = 1.875
$row44[0,19] = -1
$row44[0,20] = -1
$row44[0,21] = 2.3
$row44[0,22] = -1
$row44[0,23] = 0.875
$row44[0,24] = 0.9750000000000001
$row44[0,25] = -1
$ws.Range("E44:AD44").Value = $row44

# Row 45
$ws.Cells.Item(45, 2).Value = 6008608
$row45 = New-Object 'object[,]' 1,26
$row45[0,0] = "BodoGlimt"
$row45[0,1] = "Haugesund"
$row45[0,2] = 2
$row45[0,3] = 1
$row45[0,4] = 0
$row45[0,5] = 1
$row45[0,6] = "H"
$row45[0,7] = 1.181
$row45[0,8] = 7.5
$row45[0,9] = 11
$row45[0,10] = 1.2
$row45[0,11] = 7.5
$row45[0,12] = 11
$row45[0,13] = -2
$row45[0,14] = 1.825
$row45[0,15] = 2.025
$row45[0,16] = 3.75
$row45[0,17] = 1.95
$row45[0,18] = 1.9
$row45[0,19] = 0.2
$row45[0,20] = -1
$row45[0,21] = -1
$row45[0,22] = -1
$row45[0,23] = 1.025
$row45[0,24] = -1
$row45[0,25] = 0.8999999999999999
$ws.Range("E45:AD45").Value = $row45

# Row 51
$ws.Cells.Item(51, 2).Value = 6011405
$row51 = New-Object 'object[,]' 1,26
$row51[0,0] = "Tromso"
$row51[0,1] = "SK Brann"
$row51[0,2] = 3
$row51[0,3] = 1
$row51[0,4] = 1
$row51[0,5] = 0
$row51[0,6] = "H"
$row51[0,7] = 3.3
$row51[0,8] = 3.5
$row51[0,9] = 1.95
$row51[0,10] = 3.1
$row51[0,11] = 3.5
$row51[0,12] = 2.25
$row51[0,13] = 0.25
$row51[0,14] = 1.875
$row51[0,15] = 1.975
$row51[0,16] = 2.75
$row51[0,17] = 1.9
$row51[0,18] = 1.95
$row51[0,19] = 2.1
$row51[0,20] = -1
$row51[0,21] = -1
$row51[0,22] = 0.875
$row51[0,23] = -1
$row51[0,24] = 0.8999999999999999
$row51[0,25] = -1
$ws.Range("E51:AD51").Value = $row51

# Row 52
$ws.Cells.Item(52, 2).Value = 6011493
$row52 = New-Object 'object[,]' 1,26
$row52[0,0] = "Haugesund"
$row52[0,1] = "Lillestrom"
$row52[0,2] = 1
$row52[0,3] = 0
$row52[0,4] = 0
$row52[0,5] = 0
$row52[0,6] = "H"
$row52[0,7] = 3.3
$row52[0,8] = 3.6
$row52[0,9] = 1.95
$row52[0,10] = 3.3
$row52[0,11] = 3.6
$row52[0,12] = 2.05
$row52[0,13] = 0.25
$row52[0,14] = 2
$row52[0,15] = 1.85
$row52[0,16] = 2.75
$row52[0,17] = 1.9
$row52[0,18] = 1.95
$row52[0,19] = 2.3
$row52[0,20] = -1
$row52[0,21] = -1
$row52[0,22] = 1
$row52[0,23] = -1
$row52[0,24] = -1
$row52[0,25] = 0.95
$ws.Range("E52:AD52").Value = $row52

# Row 144
$ws.Cells.Item(144, 2).Value = 6011426
$row144 = New-Object 'object[,]' 1,26
$row144[0,0] = "SK Brann"
$row144[0,1] = "Odd BK"
$row144[0,2] = 2
$row144[0,3] = 1
$row144[0,4] = 0
$row144[0,5] = 0
$row144[0,6] = "H"
$row144[0,7] = 1.3
$row144[0,8] = 5.5
$row144[0,9] = 9
$row144[0,10] = 1.142
$row144[0,11] = 9
$row144[0,12] = 17
$row144[0,13] = -2.25
$row144[0,14] = 1.875
$row144[0,15] = 1.975
$row144[0,16] = 3.5
$row144[0,17] = 1.95
$row144[0,18] = 1.9
$row144[0,19] = 0.1419999999999999
$row144[0,20] = -1
$row144[0,21] = -1
$row144[0,22] = -1
$row144[0,23] = 0.9750000000000001
$row144[0,24] = -1
$row144[0,25] = 0.8999999999999999
$ws.Range("E144:AD144").Value = $row144

# Row 145
$ws.Cells.Item(145, 2).Value = 6011526
$row145 = New-Object 'object[,]' 1,26
$row145[0,0] = "Lillestrom"
$row145[0,1] = "Haugesund"
$row145[0,2] = 1
$row145[0,3] = 0
$row145[0,4] = 0
$row145[0,5] = 0
$row145[0,6] = "H"
$row145[0,7] = 1.7
$row145[0,8] = 4
$row145[0,9] = 4.333
$row145[0,10] = 1.8
$row145[0,11] = 3.8
$row145[0,12] = 4.2
$row145[0,13] = -0.5
$row145[0,14] = 1.86
$row145[0,15] = 2.04
$row145[0,16] = 2.5
$row145[0,17] = 1.9
$row145[0,18] = 1.95
$row145[0,19] = 0.8
$row145[0,20] = -1
$row145[0,21] = -1
$row145[0,22] = 0.8600000000000001
$row145[0,23] = -1
$row145[0,24] = -1
$row145[0,25] = 0.95
$ws.Range("E145:AD145").Value = $row145

# Row 146
$ws.Cells.Item(146, 2).Value = 6011938
$row146 = New-Object 'object[,]' 1,26
$row146[0,0] = "Sarpsborg"
$row146[0,1] = "Valerenga"
$row146[0,2] = 3
$row146[0,3] = 2
$row146[0,4] = 3
$row146[0,5] = 1
$row146[0,6] = "H"
$row146[0,7] = 1.909
$row146[0,8] = 4
$row146[0,9] = 3.5
$row146[0,10] = 1.909
$row146[0,11] = 4
$row146[0,12] = 3.5
$row146[0,13] = -0.5
$row146[0,14] = 1.9
$row146[0,15] = 2
$row146[0,16] = 3.25
$row146[0,17] = 2.025
$row146[0,18] = 1.825
$row146[0,19] = 0.909
$row146[0,20] = -1
$row146[0,21] = -1
$row146[0,22] = 0.8999999999999999
$row146[0,23] = -1
$row146[0,24] = 1.025
$row146[0,25] = -1
$ws.Range("E146:AD146").Value = $row146

# Row 152
$ws.Cells.Item(152, 2).Value = 6011428
$row152 = New-Object 'object[,]' 1,26
$row152[0,0] = "Odd BK"
$row152[0,1] = "Tromso"
$row152[0,2] = 1
$row152[0,3] = 2
$row152[0,4] = 0
$row152[0,5] = 1
$row152[0,6] = "A"
$row152[0,7] = 3.4
$row152[0,8] = 3.5
$row152[0,9] = 2.05
$row152[0,10] = 4
$row152[0,11] = 3.6
$row152[0,12] = 1.909
$row152[0,13] = 0.5
$row152[0,14] = 1.975
$row152[0,15] = 1.875
$row152[0,16] = 2.5
$row152[0,17] = 1.925
$row152[0,18] = 1.925
$row152[0,19] = -1
$row152[0,20] = -1
$row152[0,21] = 0.909
$row152[0,22] = -1
$row152[0,23] = 0.875
$row152[0,24] = 0.925
$row152[0,25] = -1
$ws.Range("E152:AD152").Value = $row152

# Row 153
$ws.Cells.Item(153, 2).Value = 6011527
$row153 = New-Object 'object[,]' 1,26
$row153[0,0] = "BodoGlimt"
$row153[0,1] = "Aalesund"
$row153[0,2] = 1
$row153[0,3] = 0
$row153[0,4] = 1
$row153[0,5] = 0
$row153[0,6] = "H"
$row153[0,7] = 1.125
$row153[0,8] = 9
$row153[0,9] = 15
$row153[0,10] = 1.1
$row153[0,11] = 10
$row153[0,12] = 21
$row153[0,13] = -2.75
$row153[0,14] = 1.88
$row153[0,15] = 2.02
$row153[0,16] = 4.25
$row153[0,17] = 1.85
$row153[0,18] = 2
$row153[0,19] = 0.1000000000000001
$row153[0,20] = -1
$row153[0,21] = -1
$row153[0,22] = -1
$row153[0,23] = 1.02
$row153[0,24] = -1
$row153[0,25] = 1
$ws.Range("E153:AD153").Value = $row153

# Row 156
$ws.Cells.Item(156, 2).Value = 6012005
$row156 = New-Object 'object[,]' 1,26
$row156[0,0] = "Valerenga"
$row156[0,1] = "Stabaek"
$row156[0,2] = 0
$row156[0,3] = 0
$row156[0,4] = 0
$row156[0,5] = 0
$row156[0,6] = "D"
$row156[0,7] = 2.1
$row156[0,8] = 3.4
$row156[0,9] = 3.4
$row156[0,10] = 1.833
$row156[0,11] = 3.75
$row156[0,12] = 4.2
$row156[0,13] = -0.5
$row156[0,14] = 1.84
$row156[0,15] = 2.06
$row156[0,16] = 2.75
$row156[0,17] = 2.1
$row156[0,18] = 1.775
$row156[0,19] = -1
$row156[0,20] = 2.75
$row156[0,21] = -1
$row156[0,22] = -1
$row156[0,23] = 1.06
$row156[0,24] = -1
$row156[0,25] = 0.7749999999999999
$ws.Range("E156:AD156").Value = $row156

# Row 157
$ws.Cells.Item(157, 2).Value = 6011939
$row157 = New-Object 'object[,]' 1,26
$row157[0,0] = "Viking FK"
$row157[0,1] = "Sarpsborg"
$row157[0,2] = 2
$row157[0,3] = 1
$row157[0,4] = 0
$row157[0,5] = 1
$row157[0,6] = "H"
$row157[0,7] = 1.75
$row157[0,8] = 4.333
$row157[0,9] = 3.8
$row157[0,10] = 1.833
$row157[0,11] = 4.5
$row157[0,12] = 3.6
$row157[0,13] = -0.5
$row157[0,14] = 1.8
$row157[0,15] = 2
$row157[0,16] = 3.75
$row157[0,17] = 1.925
$row157[0,18] = 1.925
$row157[0,19] = 0.833
$row157[0,20] = -1
$row157[0,21] = -1
$row157[0,22] = 0.8
$row157[0,23] = -1
$row157[0,24] = -1
$row157[0,25] = 0.925
$ws.Range("E157:AD157").Value = $row157

# Row 162
$ws.Cells.Item(162, 2).Value = 6011532
$row162 = New-Object 'object[,]' 1,26
$row162[0,0] = "HamKam"
$row162[0,1] = "Valerenga"
$row162[0,2] = 0
$row162[0,3] = 2
$row162[0,4] = 0
$row162[0,5] = 0
$row162[0,6] = "A"
$row162[0,7] = 3.2
$row162[0,8] = 3.6
$row162[0,9] = 2.1
$row162[0,10] = 4.2
$row162[0,11] = 3.6
$row162[0,12] = 1.85
$row162[0,13] = 0.5
$row162[0,14] = 1.975
$row162[0,15] = 1.875
$row162[0,16] = 2.5
$row162[0,17] = 1.9
$row162[0,18] = 1.95
$row162[0,19] = -1
$row162[0,20] = -1
$row162[0,21] = 0.8500000000000001
$row162[0,22] = -1
$row162[0,23] = 0.875
$row162[0,24] = -1
$row162[0,25] = 0.95
$ws.Range("E162:AD162").Value = $row162

# Row 163
$ws.Cells.Item(163, 2).Value = 6011533
$row163 = New-Object 'object[,]' 1,26
$row163[0,0] = "Stabaek"
$row163[0,1] = "Sandefjord"
$row163[0,2] = 2
$row163[0,3] = 1
$row163[0,4] = 1
$row163[0,5] = 0
$row163[0,6] = "H"
$row163[0,7] = 2
$row163[0,8] = 3.75
$row163[0,9] = 3.4
$row163[0,10] = 1.85
$row163[0,11] = 4
$row163[0,12] = 3.8
$row163[0,13] = -0.5
$row163[0,14] = 1.91
$row163[0,15] = 1.99
$row163[0,16] = 2.75
$row163[0,17] = 1.875
$row163[0,18] = 1.975
$row163[0,19] = 0.8500000000000001
$row163[0,20] = -1
$row163[0,21] = -1
$row163[0,22] = 0.9099999999999999
$row163[0,23] = -1
$row163[0,24] = 0.4375
$row163[0,25] = -0.5
$ws.Range("E163:AD163").Value = $row163

# Row 164
$ws.Cells.Item(164, 2).Value = 6012007
$row164 = New-Object 'object[,]' 1,26
$row164[0,0] = "Rosenborg"
$row164[0,1] = "Stromsgodset"
$row164[0,2] = 1
$row164[0,3] = 3
$row164[0,4] = 1
$row164[0,5] = 2
$row164[0,6] = "A"
$row164[0,7] = 1.727
$row164[0,8] = 4
$row164[0,9] = 4.2
$row164[0,10] = 1.909
$row164[0,11] = 3.8
$row164[0,12] = 3.8
$row164[0,13] = -0.5
$row164[0,14] = 1.93
$row164[0,15] = 1.97
$row164[0,16] = 3
$row164[0,17] = 2.025
$row164[0,18] = 1.825
$row164[0,19] = -1
$row164[0,20] = -1
$row164[0,21] = 2.8
$row164[0,22] = -1
$row164[0,23] = 0.97
$row164[0,24] = 1.025
$row164[0,25] = -1
$ws.Range("E164:AD164").Value = $row164

# Row 165
$ws.Cells.Item(165, 2).Value = 6011531
$row165 = New-Object 'object[,]' 1,26
$row165[0,0] = "Aalesund"
$row165[0,1] = "Viking FK"
$row165[0,2] = 0
$row165[0,3] = 4
$row165[0,4] = 0
$row165[0,5] = 2
$row165[0,6] = "A"
$row165[0,7] = 4.75
$row165[0,8] = 4.5
$row165[0,9] = 1.571
$row165[0,10] = 3.6
$row165[0,11] = 4.2
$row165[0,12] = 1.85
$row165[0,13] = 0.5
$row165[0,14] = 2.07
$row165[0,15] = 1.83
$row165[0,16] = 3.25
$row165[0,17] = 1.95
$row165[0,18] = 1.9
$row165[0,19] = -1
$row165[0,20] = -1
$row165[0,21] = 0.8500000000000001
$row165[0,22] = -1
$row165[0,23] = 0.8300000000000001
$row165[0,24] = 0.95
$row165[0,25] = -1
$ws.Range("E165:AD165").Value = $row165

# Row 166
$ws.Cells.Item(166, 2).Value = 6012006
$row166 = New-Object 'object[,]' 1,26
$row166[0,0] = "SK Brann"
$row166[0,1] = "BodoGlimt"
$row166[0,2] = 4
$row166[0,3] = 2
$row166[0,4] = 3
$row166[0,5] = 0
$row166[0,6] = "H"
$row166[0,7] = 2.15
$row166[0,8] = 4
$row166[0,9] = 2.9
$row166[0,10] = 1.95
$row166[0,11] = 4.2
$row166[0,12] = 3.4
$row166[0,13] = -0.5
$row166[0,14] = 1.9
$row166[0,15] = 1.95
$row166[0,16] = 3.5
$row166[0,17] = 1.975
$row166[0,18] = 1.875
$row166[0,19] = 0.95
$row166[0,20] = -1
$row166[0,21] = -1
$row166[0,22] = 0.8999999999999999
$row166[0,23] = -1
$row166[0,24] = 0.9750000000000001
$row166[0,25] = -1
$ws.Range("E166:AD166").Value = $row166

# Row 177
$ws.Cells.Item(177, 2).Value = 7617318
$row177 = New-Object 'object[,]' 1,26
$row177[0,0] = "Lillestrom"
$row177[0,1] = "Kristiansund BK"
$row177[0,2] = 2
$row177[0,3] = 3
$row177[0,4] = 0
$row177[0,5] = 1
$row177[0,6] = "A"
$row177[0,7] = 1.6
$row177[0,8] = 4.2
$row177[0,9] = 5
$row177[0,10] = 1.5
$row177[0,11] = 4
$row177[0,12] = 6.5
$row177[0,13] = -1
$row177[0,14] = 1.89
$row177[0,15] = 2.01
$row177[0,16] = 2.75
$row177[0,17] = 2.025
$row177[0,18] = 1.825
$row177[0,19] = -1
$row177[0,20] = -1
$row177[0,21] = 5.5
$row177[0,22] = -1
$row177[0,23] = 1.01
$row177[0,24] = 1.025
$row177[0,25] = -1
$ws.Range("E177:AD177").Value = $row177

# Row 178
$ws.Cells.Item(178, 2).Value = 7617319
$row178 = New-Object 'object[,]' 1,26
$row178[0,0] = "Molde"
$row178[0,1] = "Stromsgodset"
$row178[0,2] = 4
$row178[0,3] = 0
$row178[0,4] = 1
$row178[0,5] = 0
$row178[0,6] = "H"
$row178[0,7] = 1.45
$row178[0,8] = 4.75
$row178[0,9] = 6
$row178[0,10] = 1.4
$row178[0,11] = 5
$row178[0,12] = 8
$row178[0,13] = -1.25
$row178[0,14] = 1.9
$row178[0,15] = 1.95
$row178[0,16] = 2.75
$row178[0,17] = 1.8
$row178[0,18] = 2.05
$row178[0,19] = 0.3999999999999999
$row178[0,20] = -1
$row178[0,21] = -1
$row178[0,22] = 0.8999999999999999
$row178[0,23] = -1
$row178[0,24] = 0.8
$row178[0,25] = -1
$ws.Range("E178:AD178").Value = $row178

# Row 179
$ws.Cells.Item(179, 2).Value = 7617321
$row179 = New-Object 'object[,]' 1,26
$row179[0,0] = "Tromso"
$row179[0,1] = "SK Brann"
$row179[0,2] = 2
$row179[0,3] = 4
$row179[0,4] = 0
$row179[0,5] = 2
$row179[0,6] = "A"
$row179[0,7] = 3.1
$row179[0,8] = 3.75
$row179[0,9] = 2.1
$row179[0,10] = 3.6
$row179[0,11] = 3.5
$row179[0,12] = 2.05
$row179[0,13] = 0.5
$row179[0,14] = 1.8
$row179[0,15] = 2.05
$row179[0,16] = 2.5
$row179[0,17] = 2
$row179[0,18] = 1.85
$row179[0,19] = -1
$row179[0,20] = -1
$row179[0,21] = 1.05
$row179[0,22] = -1
$row179[0,23] = 1.05
$row179[0,24] = 1
$row179[0,25] = -1
$ws.Range("E179:AD179").Value = $row179

# Row 210
$ws.Cells.Item(210, 2).Value = 7617346
$row210 = New-Object 'object[,]' 1,26
$row210[0,0] = "KFUM"
$row210[0,1] = "SK Brann"
$row210[0,2] = 0
$row210[0,3] = 0
$row210[0,4] = 0
$row210[0,5] = 0
$row210[0,6] = "D"
$row210[0,7] = 4.2
$row210[0,8] = 4.1
$row210[0,9] = 1.727
$row210[0,10] = 4.333
$row210[0,11] = 4.333
$row210[0,12] = 1.7
$row210[0,13] = 0.75
$row210[0,14] = 1.975
$row210[0,15] = 1.875
$row210[0,16] = 2.75
$row210[0,17] = 1.825
$row210[0,18] = 2.025
$row210[0,19] = -1
$row210[0,20] = 3.333
$row210[0,21] = -1
$row210[0,22] = 0.9750000000000001
$row210[0,23] = -1
$row210[0,24] = -1
$row210[0,25] = 1.025
$ws.Range("E210:AD210").Value = $row210

# Row 211
$ws.Cells.Item(211, 2).Value = 7617345
$row211 = New-Object 'object[,]' 1,26
$row211[0,0] = "Fredrikstad"
$row211[0,1] = "Sandefjord"
$row211[0,2] = 1
$row211[0,3] = 0
$row211[0,4] = 1
$row211[0,5] = 0
$row211[0,6] = "H"
$row211[0,7] = 2
$row211[0,8] = 3.6
$row211[0,9] = 3.5
$row211[0,10] = 1.85
$row211[0,11] = 3.4
$row211[0,12] = 4.333
$row211[0,13] = -0.5
$row211[0,14] = 1.88
$row211[0,15] = 2.02
$row211[0,16] = 2.25
$row211[0,17] = 1.85
$row211[0,18] = 2
$row211[0,19] = 0.8500000000000001
$row211[0,20] = -1
$row211[0,21] = -1
$row211[0,22] = 0.8799999999999999
$row211[0,23] = -1
$row211[0,24] = -1
$row211[0,25] = 1
$ws.Range("E211:AD211").Value = $row211

# Row 212
$ws.Cells.Item(212, 2).Value = 7617347
$row212 = New-Object 'object[,]' 1,26
$row212[0,0] = "Lillestrom"
$row212[0,1] = "HamKam"
$row212[0,2] = 1
$row212[0,3] = 1
$row212[0,4] = 1
$row212[0,5] = 0
$row212[0,6] = "D"
$row212[0,7] = 1.615
$row212[0,8] = 4.1
$row212[0,9] = 5
$row212[0,10] = 1.533
$row212[0,11] = 4.2
$row212[0,12] = 5.75
$row212[0,13] = -1
$row212[0,14] = 1.925
$row212[0,15] = 1.925
$row212[0,16] = 2.75
$row212[0,17] = 1.95
$row212[0,18] = 1.9
$row212[0,19] = -1
$row212[0,20] = 3.2
$row212[0,21] = -1
$row212[0,22] = -1
$row212[0,23] = 0.925
$row212[0,24] = -1
$row212[0,25] = 0.8999999999999999
$ws.Range("E212:AD212").Value = $row212

# Row 240
$ws.Cells.Item(240, 2).Value = 7617375
$row240 = New-Object 'object[,]' 1,26
$row240[0,0] = "Lillestrom"
$row240[0,1] = "Fredrikstad"
$row240[0,2] = 0
$row240[0,3] = 3
$row240[0,4] = 0
$row240[0,5] = 1
$row240[0,6] = "A"
$row240[0,7] = 2.25
$row240[0,8] = 3.3
$row240[0,9] = 3.1
$row240[0,10] = 2.05
$row240[0,11] = 3.4
$row240[0,12] = 3.7
$row240[0,13] = -0.5
$row240[0,14] = 1.99
$row240[0,15] = 1.91
$row240[0,16] = 2.5
$row240[0,17] = 1.9
$row240[0,18] = 1.95
$row240[0,19] = -1
$row240[0,20] = -1
$row240[0,21] = 2.7
$row240[0,22] = -1
$row240[0,23] = 0.9099999999999999
$row240[0,24] = 0.8999999999999999
$row240[0,25] = -1
$ws.Range("E240:AD240").Value = $row240

# Row 241
$ws.Cells.Item(241, 2).Value = 7617376
$row241 = New-Object 'object[,]' 1,26
$row241[0,0] = "Molde"
$row241[0,1] = "Sarpsborg"
$row241[0,2] = 2
$row241[0,3] = 4
$row241[0,4] = 1
$row241[0,5] = 2
$row241[0,6] = "A"
$row241[0,7] = 1.533
$row241[0,8] = 4.75
$row241[0,9] = 5.25
$row241[0,10] = 1.38
$row241[0,11] = 5.5
$row241[0,12] = 6.25
$row241[0,13] = -1.5
$row241[0,14] = 2.03
$row241[0,15] = 1.87
$row241[0,16] = 3.25
$row241[0,17] = 1.85
$row241[0,18] = 2
$row241[0,19] = -1
$row241[0,20] = -1
$row241[0,21] = 5.25
$row241[0,22] = -1
$row241[0,23] = 0.8700000000000001
$row241[0,24] = 0.8500000000000001
$row241[0,25] = -1
$ws.Range("E241:AD241").Value = $row241

# Row 242
$ws.Cells.Item(242, 2).Value = 7617377
$row242 = New-Object 'object[,]' 1,26
$row242[0,0] = "Sandefjord"
$row242[0,1] = "Viking FK"
$row242[0,2] = 0
$row242[0,3] = 3
$row242[0,4] = 0
$row242[0,5] = 1
$row242[0,6] = "A"
$row242[0,7] = 3
$row242[0,8] = 3.5
$row242[0,9] = 2.25
$row242[0,10] = 2.8
$row242[0,11] = 3.6
$row242[0,12] = 2.3
$row242[0,13] = 0.25
$row242[0,14] = 1.8
$row242[0,15] = 2.05
$row242[0,16] = 3.25
$row242[0,17] = 1.875
$row242[0,18] = 1.975
$row242[0,19] = -1
$row242[0,20] = -1
$row242[0,21] = 1.3
$row242[0,22] = -1
$row242[0,23] = 1.05
$row242[0,24] = -0.5
$row242[0,25] = 0.4875
$ws.Range("E242:AD242").Value = $row242

# Row 243
$ws.Cells.Item(243, 2).Value = 7617374
$row243 = New-Object 'object[,]' 1,26
$row243[0,0] = "Kristiansund BK"
$row243[0,1] = "SK Brann"
$row243[0,2] = 2
$row243[0,3] = 2
$row243[0,4] = 1
$row243[0,5] = 1
$row243[0,6] = "D"
$row243[0,7] = 3.9
$row243[0,8] = 3.4
$row243[0,9] = 1.95
$row243[0,10] = 7
$row243[0,11] = 4.5
$row243[0,12] = 1.42
$row243[0,13] = 1.25
$row243[0,14] = 1.975
$row243[0,15] = 1.875
$row243[0,16] = 3.25
$row243[0,17] = 1.95
$row243[0,18] = 1.9
$row243[0,19] = -1
$row243[0,20] = 3.5
$row243[0,21] = -1
$row243[0,22] = 0.9750000000000001
$row243[0,23] = -1
$row243[0,24] = 0.95
$row243[0,25] = -1
$ws.Range("E243:AD243").Value = $row243

# Row 244
$ws.Cells.Item(244, 2).Value = 7617378
$row244 = New-Object 'object[,]' 1,26
$row244[0,0] = "Stromsgodset"
$row244[0,1] = "Haugesund"
$row244[0,2] = 2
$row244[0,3] = 0
$row244[0,4] = 1
$row244[0,5] = 0
$row244[0,6] = "H"
$row244[0,7] = 1.65
$row244[0,8] = 3.7
$row244[0,9] = 5.25
$row244[0,10] = 1.95
$row244[0,11] = 3.6
$row244[0,12] = 3.75
$row244[0,13] = -0.5
$row244[0,14] = 1.975
$row244[0,15] = 1.875
$row244[0,16] = 2.75
$row244[0,17] = 1.975
$row244[0,18] = 1.875
$row244[0,19] = 0.95
$row244[0,20] = -1
$row244[0,21] = -1
$row244[0,22] = 0.9750000000000001
$row244[0,23] = -1
$row244[0,24] = -1
$row244[0,25] = 0.875
$ws.Range("E244:AD244").Value = $row244

# Row 257 (odds update for upcoming match)
$ws.Cells.Item(257, 19).Value = 1.98
$ws.Cells.Item(257, 20).Value = 1.92
$ws.Cells.Item(257, 22).Value = 1.95
$ws.Cells.Item(257, 23).Value = 1.9
